$d = $word.ActiveDocument

# 1. Change "Enterprise " to "Enterprise Architect" in the body bullet list item.
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute("Enterprise ", $true, $false, $false, $false, $false, $true, 1, $false, "Enterprise Architect", 2)

# 2. Update cached PAGE field result in the footer (second/default footer, table cell 3)
#    from "4" to "5" to reflect the new pagination.
foreach ($sec in $d.Sections) {
    foreach ($ftrType in 1, 2, 3) {
        $ftr = $sec.Footers.Item($ftrType)
        if ($ftr.Exists) {
            $fr = $ftr.Range
            $fr.Find.ClearFormatting()
            $fr.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2)
        }
    }
}
